# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (year) / "Valor Mora" (overdue amount) rows for
# 2005 and 2003 were swapped during the database refresh; the 2004 row
# stayed the same.
#
#   Row 16 (was 2005 / 50.000)  -> 2003 / 300.000
#   Row 17 (2004 / 300.000)     -> unchanged
#   Row 18 (was 2003 / 300.000) -> 2005 / 50.000

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2003"
$ws.Range("F16").Value = 300000

$ws.Range("E18").Value = "2005"
$ws.Range("F18").Value = 50000
